$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")
$ws.Range("A1").Value = $ws.Range("A1").Value
